$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.536.07"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "3.253.22"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'577.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "'181.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.29%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "3.252.49"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("E10").Value = "  +5.53%  "
$ws.Range("D11").Value = "'6.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +5.20%  "
$ws.Range("D13").Value = "3.817.02"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "'28.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").Value = "67.522.15"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "3.257.82"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "'5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.42%  "
$ws.Range("D21").Value = "'376.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.70%  "
$ws.Range("E22").Value = "  +5.41%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").Value = "'0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").Value = "'9.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  +7.60%  "
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.35%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'6.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("D36").Value = "'163.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.48%  "
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").Value = "'0.852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("D40").Value = "'26.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "'6.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.75%  "
$ws.Range("D42").Value = "'4.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.45%  "
$ws.Range("D43").Value = "'2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.21%  "
$ws.Range("D44").Value = "'362.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.81%  "
$ws.Range("D45").Value = "2.735.35"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").Value = "'25.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.42%  "
$ws.Range("D47").Value = "'40.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("D48").Value = "'0.0681"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.58%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("E50").Value = "  +6.95%  "
$ws.Range("E51").Value = "  +0.48%  "
